$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$arr = New-Object 'object[,]' 25,9
$arr[0,0] = "model_9_6_24"
$arr[0,1] = -0.1283255146355962
$arr[0,2] = -1.910858872220313
$arr[0,3] = -10.89708812026705
$arr[0,4] = -2.878134316471546
$arr[0,5] = 1.248723864555359
$arr[0,6] = 3.796500682830811
$arr[0,7] = 4.43040132522583
$arr[0,8] = 4.094807624816895
$arr[1,0] = "model_9_6_23"
$arr[1,1] = -0.1175926631363253
$arr[1,2] = -1.878948554820687
$arr[1,3] = -10.79879026942878
$arr[1,4] = -2.840949312628336
$arr[1,5] = 1.236845850944519
$arr[1,6] = 3.754881620407104
$arr[1,7] = 4.393795490264893
$arr[1,8] = 4.055545806884766
$arr[2,0] = "model_9_6_22"
$arr[2,1] = -0.1161484145237721
$arr[2,2] = -1.873616565857722
$arr[2,3] = -10.78728694567486
$arr[2,4] = -2.835554898766522
$arr[2,5] = 1.235247373580933
$arr[2,6] = 3.747926950454712
$arr[2,7] = 4.389511585235596
$arr[2,8] = 4.049849510192871
$arr[3,0] = "model_9_6_21"
$arr[3,1] = -0.1001099567280483
$arr[3,2] = -1.87791214129056
$arr[3,3] = -10.41556299805001
$arr[3,4] = -2.77666804867948
$arr[3,5] = 1.217497587203979
$arr[3,6] = 3.75352954864502
$arr[3,7] = 4.251084327697754
$arr[3,8] = 3.987672805786133
$arr[4,0] = "model_9_6_20"
$arr[4,1] = -0.08049776147066279
$arr[4,2] = -1.841060330504713
$arr[4,3] = -10.13398808235953
$arr[4,4] = -2.705835971226788
$arr[4,5] = 1.195792555809021
$arr[4,6] = 3.705465078353882
$arr[4,7] = 4.146227359771729
$arr[4,8] = 3.912883043289185
$arr[5,0] = "model_9_6_18"
$arr[5,1] = -0.06467367812617608
$arr[5,2] = -1.767760024676936
$arr[5,3] = -10.05160282363792
$arr[5,4] = -2.64422755500111
$arr[5,5] = 1.178279995918274
$arr[5,6] = 3.60986328125
$arr[5,7] = 4.115547657012939
$arr[5,8] = 3.847832679748535
$arr[6,0] = "model_9_6_19"
$arr[6,1] = -0.05942387465466248
$arr[6,2] = -1.75533445498699
$arr[6,3] = -9.998962426948861
$arr[6,4] = -2.627364559908933
$arr[6,5] = 1.172470092773438
$arr[6,6] = 3.593657255172729
$arr[6,7] = 4.095944881439209
$arr[6,8] = 3.830027341842651
$arr[7,0] = "model_9_6_17"
$arr[7,1] = -0.055183459303358
$arr[7,2] = -1.744990195451381
$arr[7,3] = -9.943828209613248
$arr[7,4] = -2.611447205940738
$arr[7,5] = 1.167777180671692
$arr[7,6] = 3.580165386199951
$arr[7,7] = 4.075413227081299
$arr[7,8] = 3.813220739364624
$arr[8,0] = "model_9_6_2"
$arr[8,1] = -0.01146815703426518
$arr[8,2] = -1.335638823397878
$arr[8,3] = -7.694564033629213
$arr[8,4] = -1.970440701455991
$arr[8,5] = 1.119397163391113
$arr[8,6] = 3.046266794204712
$arr[8,7] = 3.237801313400269
$arr[8,8] = 3.136400938034058
$arr[9,0] = "model_9_6_16"
$arr[9,1] = 0.03498124306223693
$arr[9,2] = -1.555555065195551
$arr[9,3] = -8.80460496607221
$arr[9,4] = -2.29848993026276
$arr[9,5] = 1.067991375923157
$arr[9,6] = 3.333094120025635
$arr[9,7] = 3.65117359161377
$arr[9,8] = 3.482778549194336
$arr[10,0] = "model_9_6_15"
$arr[10,1] = 0.07428562569929587
$arr[10,2] = -1.468377596766559
$arr[10,3] = -8.278693297962988
$arr[10,4] = -2.154193867266348
$arr[10,5] = 1.024492979049683
$arr[10,6] = 3.2193922996521
$arr[10,7] = 3.455327272415161
$arr[10,8] = 3.33042049407959
$arr[11,0] = "model_9_6_14"
$arr[11,1] = 0.08223093416710925
$arr[11,2] = -1.456444002023354
$arr[11,3] = -8.116030044379555
$arr[11,4] = -2.119392629823556
$arr[11,5] = 1.015699982643127
$arr[11,6] = 3.203827619552612
$arr[11,7] = 3.394752502441406
$arr[11,8] = 3.293674945831299
$arr[12,0] = "model_9_6_13"
$arr[12,1] = 0.148237314173244
$arr[12,2] = -1.303569302870291
$arr[12,3] = -7.265257234822684
$arr[12,4] = -1.878216437065231
$arr[12,5] = 0.942650318145752
$arr[12,6] = 3.004440069198608
$arr[12,7] = 3.077930212020874
$arr[12,8] = 3.039024114608765
$arr[13,0] = "model_9_6_12"
$arr[13,1] = 0.1919947603557329
$arr[13,2] = -1.137688971457123
$arr[13,3] = -6.923094362682461
$arr[13,4] = -1.712948878186062
$arr[13,5] = 0.8942237496376038
$arr[13,6] = 2.788089990615845
$arr[13,7] = 2.95051097869873
$arr[13,8] = 2.864522933959961
$arr[14,0] = "model_9_6_11"
$arr[14,1] = 0.2022856442654577
$arr[14,2] = -1.110979519572607
$arr[14,3] = -6.80781669149287
$arr[14,4] = -1.67634817050734
$arr[14,5] = 0.8828348517417908
$arr[14,6] = 2.753254413604736
$arr[14,7] = 2.90758228302002
$arr[14,8] = 2.82587718963623
$arr[15,0] = "model_9_6_10"
$arr[15,1] = 0.2109969569673066
$arr[15,2] = -1.088653773865948
$arr[15,3] = -6.705064699902424
$arr[15,4] = -1.644696694930119
$arr[15,5] = 0.8731938600540161
$arr[15,6] = 2.724135637283325
$arr[15,7] = 2.869317770004272
$arr[15,8] = 2.792457342147827
$arr[16,0] = "model_9_6_3"
$arr[16,1] = 0.2254262492117556
$arr[16,2] = -1.006675850214173
$arr[16,3] = -5.684244412634976
$arr[16,4] = -1.421664300615271
$arr[16,5] = 0.8572248816490173
$arr[16,6] = 2.617215633392334
$arr[16,7] = 2.489170789718628
$arr[16,8] = 2.556964159011841
$arr[17,0] = "model_9_6_9"
$arr[17,1] = 0.2378867930180456
$arr[17,2] = -1.033303568224882
$arr[17,3] = -6.185898919104341
$arr[17,4] = -1.522333608185216
$arr[17,5] = 0.843434751033783
$arr[17,6] = 2.651944875717163
$arr[17,7] = 2.675983905792236
$arr[17,8] = 2.663257837295532
$arr[18,0] = "model_9_6_4"
$arr[18,1] = 0.2477079273764463
$arr[18,2] = -0.953950576035673
$arr[18,3] = -5.676655449155457
$arr[18,4] = -1.385920989355253
$arr[18,5] = 0.8325656652450562
$arr[18,6] = 2.54844856262207
$arr[18,7] = 2.486344814300537
$arr[18,8] = 2.519223928451538
$arr[19,0] = "model_9_6_6"
$arr[19,1] = 0.2506069678052927
$arr[19,2] = -1.025611333464981
$arr[19,3] = -5.721743190999925
$arr[19,4] = -1.440264410424478
$arr[19,5] = 0.8293572664260864
$arr[19,6] = 2.641912221908569
$arr[19,7] = 2.503135204315186
$arr[19,8] = 2.576603412628174
$arr[20,0] = "model_9_6_7"
$arr[20,1] = 0.2536279814114499
$arr[20,2] = -1.013268202835603
$arr[20,3] = -5.744644925929951
$arr[20,4] = -1.435994831937101
$arr[20,5] = 0.8260138630867004
$arr[20,6] = 2.625813961029053
$arr[20,7] = 2.511663675308228
$arr[20,8] = 2.572095394134521
$arr[21,0] = "model_9_6_8"
$arr[21,1] = 0.2602334689391798
$arr[21,2] = -0.9825700866144043
$arr[21,3] = -5.769112052067409
$arr[21,4] = -1.419981878324873
$arr[21,5] = 0.8187035322189331
$arr[21,6] = 2.585775375366211
$arr[21,7] = 2.520774841308594
$arr[21,8] = 2.555187702178955
$arr[22,0] = "model_9_6_5"
$arr[22,1] = 0.270879167659998
$arr[22,2] = -0.9387075280484203
$arr[22,3] = -5.641227338695455
$arr[22,4] = -1.370071646502482
$arr[22,5] = 0.8069219589233398
$arr[22,6] = 2.528567790985107
$arr[22,7] = 2.473151445388794
$arr[22,8] = 2.50248908996582
$arr[23,0] = "model_9_6_1"
$arr[23,1] = 0.3754979927268688
$arr[23,2] = 0.4016680086832379
$arr[23,3] = -3.499014636768559
$arr[23,4] = -0.1379872962511344
$arr[23,5] = 0.6911396980285645
$arr[23,6] = 0.7803771495819092
$arr[23,7] = 1.675404906272888
$arr[23,8] = 1.20156729221344
$arr[24,0] = "model_9_6_0"
$arr[24,1] = 0.4409215921293985
$arr[24,2] = 0.6535744171904152
$arr[24,3] = -2.920473254677894
$arr[24,4] = 0.122768615503896
$arr[24,5] = 0.6187350153923035
$arr[24,6] = 0.4518270492553711
$arr[24,7] = 1.459959626197815
$arr[24,8] = 0.9262428283691406
$ws.Range("A2:I26").Value = $arr
